# Apply the "Talk" -> "Direct instruction/Direct Instruction" label change
# as described in the commit message / diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Abb. Workshop - Talk" -> "Abb. Workshop - Direct instruction"
$ws.Range("A3").Value = "Abb. Workshop - Direct instruction"

# "Workshop - Talk" -> "Workshop - Direct Instruction"
$ws.Range("A4").Value = "Workshop - Direct Instruction"

# "Abb. Workshop - Talk" -> "Abb. Workshop - Direct instruction"
$ws.Range("A6").Value = "Abb. Workshop - Direct instruction"

# "Workshop - Talk" -> "Workshop - Direct Instruction"
$ws.Range("A7").Value = "Workshop - Direct Instruction"
